# React Session 1.pptx - "useEffect hook" slide (slide 12)
#
# The SmartArt list diagram's last bullet explained the dependency array
# but trailed off ("...if left empty then"). Update it to finish the
# thought: an empty dependency array means useEffect only runs once, on
# mount ("at the start").

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(12)

$oldText = "Dependency array, if left empty then"
$newText = "Dependency array, if left as empty array, useEffect is only called at the start"

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.HasSmartArt) {
        $sa = $sh.SmartArt
        for ($j = 1; $j -le $sa.AllNodes.Count; $j++) {
            $node = $sa.AllNodes.Item($j)
            $tr = $node.TextFrame2.TextRange
            if ($tr.Text -eq $oldText) {
                $tr.Text = $newText
            }
        }
    }
}
